$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Add EXI Devices Loop A")
$ws2 = $wb.Worksheets.Item("Add IS Devices to EXI800")

# --- Sheet view scroll position changes ---
$ws1.Application.ActiveWindow.ScrollColumn = 2
$ws2.Application.ActiveWindow.ScrollColumn = 2

# --- Sheet 2: comment / debug text cells ---

# I4: "//Verify IS Unit after deleting IS device" with special (grayish) comment font color
$ws2.Range("I4").ClearFormats()
$ws2.Range("I4").Value = "//Verify IS Unit after deleting IS device"
$ws2.Range("I4").Font.Color = 8756119

# H5: new numeric value 41
$ws2.Range("H5").Value = 41

# I5: the long rich-text code comment line
$codeText = "  line number462 need to change to j,8          sISUnits = ((Range)Excel_Utilities.ExcelRange.Cells[j+1,8]).Value.ToString();"
$ws2.Range("I5").ClearFormats()
$ws2.Range("I5").Value = $codeText

$ws2.Range("I5").Characters(103, 1).Font.Color = 9109504
$ws2.Range("I5").Characters(105, 1).Font.Color = 9109504
$ws2.Range("I5").Characters(115, 8).Font.Bold = $true
$ws2.Range("I5").Characters(115, 8).Font.Color = 7346457

# H6: new numeric value 43
$ws2.Range("H6").Value = 43

# I6: a string of spaces
$ws2.Range("I6").ClearFormats()
$ws2.Range("I6").Value = "            "

# --- Selection adjustments to match final view state ---
$ws2.Range("H8").Select()
$ws1.Range("I8").Select()
$ws2.Select()
$ws2.Range("H8").Select()
